# Updates the crypto price/volume table (rows 2-51) on the active sheet
# to match the refreshed feed values from the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '34.406.88'
$ws.Range("E2").Value = '  +0.95%  '

# Row 3
$ws.Range("D3").Value = '1.796.62'
$ws.Range("E3").Value = '  +0.68%  '

# Row 4
$ws.Range("E4").Value = '  -0.26%  '

# Row 5
$ws.Range("D5").Value = '''226.70'
$ws.Range("E5").Value = '  +0.35%  '

# Row 6
$ws.Range("E6").Value = '  +1.69%  '

# Row 7
$ws.Range("E7").Value = '  -0.26%  '

# Row 8
$ws.Range("E8").Value = '  +1.91%  '

# Row 9
$ws.Range("E9").Value = '  +1.73%  '

# Row 10
$ws.Range("E10").Value = '  +0.82%  '

# Row 11
$ws.Range("D11").Value = '''0.0949'
$ws.Range("E11").Value = '  +0.59%  '

# Row 12
$ws.Range("D12").Value = '2.055.04'
$ws.Range("E12").Value = '  +0.65%  '

# Row 13
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").Value = '''11.05'
$ws.Range("E13").Value = '  -1.22%  '

# Row 14
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.793.60'
$ws.Range("E14").Value = '  +0.42%  '

# Row 15
$ws.Range("E15").Value = '  +2.13%  '

# Row 16
$ws.Range("D16").Value = '34.356.83'
$ws.Range("E16").Value = '  +0.91%  '

# Row 17
$ws.Range("E17").Value = '  +1.21%  '

# Row 18
$ws.Range("D18").Value = '''68.35'
$ws.Range("E18").Value = '  +0.77%  '

# Row 19
$ws.Range("D19").Value = '0.0₃0802'
$ws.Range("E19").Value = '  +3.32%  '

# Row 20
$ws.Range("D20").Value = '''246.90'
$ws.Range("E20").Value = '  +0.69%  '

# Row 21
$ws.Range("D21").Value = '''10.99'
$ws.Range("E21").Value = '  +2.01%  '

# Row 22
$ws.Range("E22").Value = '  -0.16%  '

# Row 23
$ws.Range("E23").Value = '  +2.17%  '

# Row 24
$ws.Range("E24").Value = '  +1.11%  '

# Row 25
$ws.Range("D25").Value = '''162.38'
$ws.Range("E25").Value = '  +0.77%  '

# Row 26
$ws.Range("D26").Value = '''7.20'
$ws.Range("E26").Value = '  +1.03%  '

# Row 27
$ws.Range("D27").Value = '''16.41'
$ws.Range("E27").Value = '  +0.87%  '

# Row 28
$ws.Range("E28").Value = '  +2.40%  '

# Row 30
$ws.Range("D30").Value = '''3.93'
$ws.Range("E30").Value = '  +9.39%  '

# Row 31
$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").Value = '''0.0522'
$ws.Range("E31").Value = '  +1.06%  '

# Row 32
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = '''1.23'
$ws.Range("E32").Value = '  +0.20%  '

# Row 33
$ws.Range("D33").Value = '''3.78'
$ws.Range("E33").Value = '  +3.84%  '

# Row 34
$ws.Range("E34").Value = '  +1.26%  '

# Row 35
$ws.Range("D35").Value = '1.442.32'
$ws.Range("E35").Value = '  -0.73%  '

# Row 36
$ws.Range("D36").Value = '''2.61'
$ws.Range("E36").Value = '  +8.92%  '

# Row 37
$ws.Range("E37").Value = '  +3.22%  '

# Row 38
$ws.Range("D38").Value = '''1.06'
$ws.Range("E38").Value = '  +1.90%  '

# Row 39
$ws.Range("E39").Value = '  -0.86%  '

# Row 40
$ws.Range("D40").Value = '''83.27'
$ws.Range("E40").Value = '  +4.46%  '

# Row 41
$ws.Range("E41").Value = '  +1.02%  '

# Row 42
$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D42").Value = '''0.934'
$ws.Range("E42").Value = '  +1.92%  '

# Row 43
$ws.Range("E43").Value = '  +2.73%  '

# Row 44
$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").Value = '''13.97'
$ws.Range("E44").Value = '  +4.61%  '

# Row 45
$ws.Range("D45").Value = '''0.0521'
$ws.Range("E45").Value = '  +2.50%  '

# Row 46
$ws.Range("D46").Value = '''6.08'
$ws.Range("E46").Value = '  +0.95%  '

# Row 47
$ws.Range("E47").Value = '  -0.27%  '

# Row 48
$ws.Range("D48").Value = '1.948.82'
$ws.Range("E48").Value = '  +0.29%  '

# Row 49
$ws.Range("D49").Value = '''105.66'
$ws.Range("E49").Value = '  -1.30%  '

# Row 50
$ws.Range("E50").Value = '  -0.25%  '

# Row 51
$ws.Range("E51").Value = '  -4.64%  '
